# Fruta / hortaliza, semanal
# Insert two new daily records (rows 156-157) into the "Macroferia Regional
# de Talca - Mandarina" sheet, pushing the previously-existing rows
# 156..242 down to 158..244.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 156 (shifts old rows 156.. downward by 2)
$ws.Range("A156:A157").EntireRow.Insert()

# New row 156
$ws.Range("A156").Value = 5
$ws.Range("B156").Value = "Macroferia Regional de Talca"
$ws.Range("C156").Value = "Maule"
$ws.Range("D156").Value = 44567
$ws.Range("E156").Value = 7
$ws.Range("F156").Value = "Fruta"
$ws.Range("G156").Value = 100102
$ws.Range("H156").Value = "Cítricos"
$ws.Range("I156").Value = 100102004
$ws.Range("J156").Value = "Mandarina"
$ws.Range("K156").Value = "Marisol"
$ws.Range("L156").Value = "Primera"
$ws.Range("M156").Value = 180
$ws.Range("N156").Value = 8000
$ws.Range("O156").Value = 8000
$ws.Range("P156").Value = 8000
$ws.Range("Q156").Value = "`$/caja 18 kilos"
$ws.Range("R156").Value = "Región de O'Higgins"
$ws.Range("S156").Value = 444
$ws.Range("T156").Value = 18

# New row 157
$ws.Range("A157").Value = 5
$ws.Range("B157").Value = "Macroferia Regional de Talca"
$ws.Range("C157").Value = "Maule"
$ws.Range("D157").Value = 44567
$ws.Range("E157").Value = 7
$ws.Range("F157").Value = "Fruta"
$ws.Range("G157").Value = 100102
$ws.Range("H157").Value = "Cítricos"
$ws.Range("I157").Value = 100102004
$ws.Range("J157").Value = "Mandarina"
$ws.Range("K157").Value = "Murcott"
$ws.Range("L157").Value = "Primera"
$ws.Range("M157").Value = 200
$ws.Range("N157").Value = 8000
$ws.Range("O157").Value = 8000
$ws.Range("P157").Value = 8000
$ws.Range("Q157").Value = "`$/caja 18 kilos"
$ws.Range("R157").Value = "Región de O'Higgins"
$ws.Range("S157").Value = 444
$ws.Range("T157").Value = 18
